$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: remove the old unit-labels row (row 2). Remaining data rows
# (old 3,4,5) shift up to become rows 2,3,4. ---
$ws.Rows.Item(2).Delete()

# --- Step 2: write the new header row (row 1) ---
# A1:E1 carry no explicit style (plain default formatting); clear any
# leftover formatting these cells may have inherited (e.g. old E1).
$ws.Range("A1:E1").ClearFormats()
$ws.Cells.Item(1, 1).Value = "idx"
$ws.Cells.Item(1, 2).Value = "idx2"
$ws.Cells.Item(1, 3).Value = "Name"
$ws.Cells.Item(1, 4).Value = "Date Start"
$ws.Cells.Item(1, 5).Value = "Date End"
$ws.Cells.Item(1, 6).Value = "(m3/s)"
$ws.Cells.Item(1, 7).Value = "(MW1)"
$ws.Cells.Item(1, 8).Value = "(MW2)"
$ws.Cells.Item(1, 9).Value = "(GWh) Winter"
$ws.Cells.Item(1, 10).Value = "(GWh) Summer"
$ws.Cells.Item(1, 11).Value = "(GWh) Year"

# --- Step 3: apply the "font only" style (Arial 9) used on F1:K1 -----
# (matches the workbook's existing body font, but via a style record
# that only flags applyFont, not applyNumberFormat)
$wb.Styles.Add("HeaderFont9")
$hs = $wb.Styles.Item("HeaderFont9")
$hs.Font.Name = "Arial"
$hs.Font.Size = 9
$ws.Range("F1:K1").Style = "HeaderFont9"
$wb.Styles.Item("HeaderFont9").Delete()

# --- Step 4: selection / view state ---
$ws.Range("A2:K2").Select()
